# tests/reports/df_filter_frame.xlsx
# Commit: "Reports: datetime and arithmetic filters / breaking change: filter
# string args must be Excel strings (#1674)"
#
# The functional change is that filter string arguments must now be quoted
# Excel strings, so the two template cells on Sheet1 that read
#   {{ df2 | maxrows(2, Other, 0)}}
#   {{ df2 | noheader | maxrows(2, Other, 0) }}
# become
#   {{ df2 | maxrows(2, "Other", 0)}}
#   {{ df2 | noheader | maxrows(2, "Other", 0) }}
# Additionally the active selection on Sheet1 moved to A7.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A3").Value = '{{ df2 | maxrows(2, "Other", 0)}}'
$ws1.Range("A6").Value = '{{ df2 | noheader | maxrows(2, "Other", 0) }}'

$ws1.Activate()
$ws1.Range("A7").Select() | Out-Null
